# Remove obsolete reseek-sensitive results
# - Update SF!G7:G42 with the new (recomputed) sensitivity values
# - Remove the now-obsolete tail of the GTalign series: G43 and rows 45:87
# - Update the active selection to I40

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SF")

$newValues = @(
    [double]"1.5339999999999999E-2",
    [double]"2.266E-2",
    [double]"3.3000000000000002E-2",
    [double]"4.6829999999999997E-2",
    [double]"6.1190000000000001E-2",
    [double]"8.1079999999999999E-2",
    [double]"0.1046",
    [double]"0.12870000000000001",
    [double]"0.15790000000000001",
    [double]"0.1946",
    [double]"0.23680000000000001",
    [double]"0.28460000000000002",
    [double]"0.34079999999999999",
    [double]"0.40810000000000002",
    [double]"0.4899",
    [double]"0.57689999999999997",
    [double]"0.68110000000000004",
    [double]"0.81289999999999996",
    [double]"0.95940000000000003",
    [double]"1.1279999999999999",
    [double]"1.3340000000000001",
    [double]"1.548",
    [double]"1.7929999999999999",
    [double]"2.08",
    [double]"2.4009999999999998",
    [double]"2.7679999999999998",
    [double]"3.1880000000000002",
    [double]"3.6619999999999999",
    [double]"4.1689999999999996",
    [double]"4.7380000000000004",
    [double]"5.3369999999999997",
    [double]"5.9909999999999997",
    [double]"6.7069999999999999",
    [double]"7.4509999999999996",
    [double]"8.266",
    [double]"9.1129999999999995"
)

$startRow = 7
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 7).Value = $newValues[$i]
}

# Clear the obsolete tail of column G (old row 43 value, 9.853) and the
# entirely-obsolete rows 45:87 that only ever held column-G data.
$ws.Range("G43").ClearContents()
$ws.Range("A45:G87").ClearContents()

$ws.Range("I40").Select()

$wb.Save()
